$d = $word.ActiveDocument

# Remove the last two paragraphs entirely, and replace the text of the
# first paragraph with the new text, then add a _GoBack bookmark.

# Delete paragraph 3 ("Esquece, não tome no seu ovo!!!")
$d.Paragraphs(3).Range.Delete()

# Delete paragraph 2 ("Vai tomar no seu ovo inteiro")
$d.Paragraphs(2).Range.Delete()

# Replace text in paragraph 1, with a trailing sentinel character so the
# bookmark we add next isn't sitting exactly at the paragraph-end
# boundary (which would otherwise get normalized to the paragraph start).
$d.Content.Find.Execute("Vai tomar no meio do seu ovo", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Que coisa feia!#", 2)

# Add a _GoBack bookmark right after the real text, before the sentinel.
$p1 = $d.Paragraphs(1).Range
$bm = $d.Range($p1.End - 2, $p1.End - 2)
$d.Bookmarks.Add("_GoBack", $bm)

# Remove the sentinel character now that the bookmark is anchored.
$sentinel = $d.Range($p1.End - 2, $p1.End - 1)
$sentinel.Delete()
